$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the changed cells to remain plain text (matching the source
# inlineStr cells) instead of being auto-coerced to numbers by Excel,
# then restore the original (default) cell style afterwards.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "61.898.38"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.410.32"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "409.92"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "128.86"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  +5.92%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +5.56%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "43.40"
$ws.Range("E11").Value = "  +2.54%  "
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "0.0000221"
$ws.Range("E12").Value = "  +30.93%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "9.35"
$ws.Range("E13").Value = "  +9.10%  "
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "21.37"
$ws.Range("E15").Value = "  +6.95%  "
$ws.Range("D16").Value = "3.952.72"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "3.426.21"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  +8.04%  "
$ws.Range("E19").Value = "  +6.19%  "
$ws.Range("D20").Value = "61.880.19"
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").Value = "449.18"
$ws.Range("E21").Value = "  +42.62%  "
$ws.Range("D22").Value = "91.54"
$ws.Range("E22").Value = "  +7.78%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "13.20"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "3.30"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").Value = "9.33"
$ws.Range("E26").Value = "  +13.83%  "
$ws.Range("D27").Value = "33.17"
$ws.Range("E27").Value = "  +10.48%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").Value = "2.72"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "12.04"
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "42.52"
$ws.Range("E34").Value = "  -4.17%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "0.0505"
$ws.Range("E36").Value = "  +3.53%  "
$ws.Range("D37").Value = "53.80"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +8.50%  "
$ws.Range("E40").Value = "  +1.58%  "
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "0.319"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "4.36"
$ws.Range("E43").Value = "  +10.72%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "143.75"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.58"
$ws.Range("E45").Value = "  +15.68%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").Value = "2.01"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "16.67"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("E48").Value = "  +22.52%  "
$ws.Range("E49").Value = "  +5.33%  "
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  +6.69%  "
$ws.Range("D51").Value = "3.753.49"
$ws.Range("E51").Value = "  -0.60%  "

$dataRange.Style = "Normal"
